$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so the cells below can be written.
$ws.Unprotect()

# Update the disclaimer date text (A9) from 2021-04-21 to 2021-04-22
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-22 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.253330973751845
$ws.Range("E2").Value = -0.01687356875979285

$ws.Range("D3").Value = 0.2474778273504167
$ws.Range("E3").Value = -0.01110478359908884

$ws.Range("D4").Value = 0.2471593941042721
$ws.Range("E4").Value = -0.005051005249083973

$ws.Range("D5").Value = 0.2520318047934663
$ws.Range("E5").Value = -0.007355838696965655

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = -0.010125094022433
